$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 4).Value = 44784
$ws.Cells.Item(2, 10).Value = 1000
$ws.Cells.Item(2, 11).Value = 1200
$ws.Cells.Item(2, 12).Value = 1300
$ws.Cells.Item(2, 13).Value = 1250
$ws.Cells.Item(2, 16).Value = 1250

# Row 3
$ws.Cells.Item(3, 4).Value = 44453
$ws.Cells.Item(3, 10).Value = 1000
$ws.Cells.Item(3, 11).Value = 800
$ws.Cells.Item(3, 12).Value = 900
$ws.Cells.Item(3, 13).Value = 850
$ws.Cells.Item(3, 16).Value = 850

# Row 4
$ws.Cells.Item(4, 4).Value = 44229
$ws.Cells.Item(4, 10).Value = 1500
$ws.Cells.Item(4, 11).Value = 1400
$ws.Cells.Item(4, 12).Value = 1500
$ws.Cells.Item(4, 13).Value = 1450
$ws.Cells.Item(4, 16).Value = 1450

# Row 5
$ws.Cells.Item(5, 4).Value = 44476
$ws.Cells.Item(5, 10).Value = 900
$ws.Cells.Item(5, 11).Value = 700
$ws.Cells.Item(5, 12).Value = 800
$ws.Cells.Item(5, 13).Value = 750
$ws.Cells.Item(5, 16).Value = 750

# Row 6
$ws.Cells.Item(6, 4).Value = 44484
$ws.Cells.Item(6, 10).Value = 900
$ws.Cells.Item(6, 11).Value = 750
$ws.Cells.Item(6, 12).Value = 800
$ws.Cells.Item(6, 13).Value = 775
$ws.Cells.Item(6, 16).Value = 775

# Row 7
$ws.Cells.Item(7, 4).Value = 44656
$ws.Cells.Item(7, 10).Value = 1000
$ws.Cells.Item(7, 11).Value = 900
$ws.Cells.Item(7, 12).Value = 1000
$ws.Cells.Item(7, 13).Value = 950
$ws.Cells.Item(7, 16).Value = 950

# Row 8
$ws.Cells.Item(8, 4).Value = 44550
$ws.Cells.Item(8, 10).Value = 1300

# Row 9
$ws.Cells.Item(9, 4).Value = 44243
$ws.Cells.Item(9, 10).Value = 1200
$ws.Cells.Item(9, 11).Value = 1200
$ws.Cells.Item(9, 12).Value = 1300
$ws.Cells.Item(9, 13).Value = 1250
$ws.Cells.Item(9, 16).Value = 1250

# Row 10
$ws.Cells.Item(10, 4).Value = 44175
$ws.Cells.Item(10, 10).Value = 1600
$ws.Cells.Item(10, 11).Value = 1000
$ws.Cells.Item(10, 12).Value = 1200
$ws.Cells.Item(10, 13).Value = 1100
$ws.Cells.Item(10, 16).Value = 1100

# Row 11
$ws.Cells.Item(11, 4).Value = 44673
$ws.Cells.Item(11, 10).Value = 900
$ws.Cells.Item(11, 11).Value = 1300
$ws.Cells.Item(11, 12).Value = 1400
$ws.Cells.Item(11, 13).Value = 1350
$ws.Cells.Item(11, 16).Value = 1350

# Row 12
$ws.Cells.Item(12, 4).Value = 44407
$ws.Cells.Item(12, 10).Value = 1000
$ws.Cells.Item(12, 11).Value = 1200
$ws.Cells.Item(12, 12).Value = 1300
$ws.Cells.Item(12, 13).Value = 1250
$ws.Cells.Item(12, 16).Value = 1250

# Row 13
$ws.Cells.Item(13, 4).Value = 44883
$ws.Cells.Item(13, 10).Value = 800
$ws.Cells.Item(13, 11).Value = 550
$ws.Cells.Item(13, 12).Value = 600
$ws.Cells.Item(13, 13).Value = 575
$ws.Cells.Item(13, 16).Value = 575

# Row 14
$ws.Cells.Item(14, 4).Value = 44638
$ws.Cells.Item(14, 10).Value = 1000
$ws.Cells.Item(14, 11).Value = 900
$ws.Cells.Item(14, 12).Value = 950
$ws.Cells.Item(14, 13).Value = 925
$ws.Cells.Item(14, 16).Value = 925

# Row 15
$ws.Cells.Item(15, 4).Value = 44607
$ws.Cells.Item(15, 10).Value = 900
$ws.Cells.Item(15, 11).Value = 1300
$ws.Cells.Item(15, 12).Value = 1400
$ws.Cells.Item(15, 13).Value = 1350
$ws.Cells.Item(15, 16).Value = 1350

# Row 16
$ws.Cells.Item(16, 4).Value = 44687
$ws.Cells.Item(16, 10).Value = 1000
$ws.Cells.Item(16, 11).Value = 1200
$ws.Cells.Item(16, 12).Value = 1300
$ws.Cells.Item(16, 13).Value = 1250
$ws.Cells.Item(16, 16).Value = 1250

# Row 18
$ws.Cells.Item(18, 4).Value = 44341
$ws.Cells.Item(18, 10).Value = 1300
$ws.Cells.Item(18, 11).Value = 900
$ws.Cells.Item(18, 12).Value = 1000
$ws.Cells.Item(18, 13).Value = 950
$ws.Cells.Item(18, 16).Value = 950

# Row 19
$ws.Cells.Item(19, 4).Value = 44284
$ws.Cells.Item(19, 10).Value = 1500
$ws.Cells.Item(19, 11).Value = 800
$ws.Cells.Item(19, 12).Value = 850
$ws.Cells.Item(19, 13).Value = 825
$ws.Cells.Item(19, 16).Value = 825

# Row 20
$ws.Cells.Item(20, 4).Value = 44442
$ws.Cells.Item(20, 10).Value = 1250
$ws.Cells.Item(20, 11).Value = 850
$ws.Cells.Item(20, 12).Value = 900
$ws.Cells.Item(20, 13).Value = 875
$ws.Cells.Item(20, 16).Value = 875

# Row 21
$ws.Cells.Item(21, 4).Value = 44291
$ws.Cells.Item(21, 10).Value = 1000
$ws.Cells.Item(21, 11).Value = 1000
$ws.Cells.Item(21, 12).Value = 1200
$ws.Cells.Item(21, 13).Value = 1100
$ws.Cells.Item(21, 16).Value = 1100

# Row 22
$ws.Cells.Item(22, 4).Value = 44449
$ws.Cells.Item(22, 10).Value = 1300
$ws.Cells.Item(22, 11).Value = 900
$ws.Cells.Item(22, 12).Value = 950
$ws.Cells.Item(22, 13).Value = 925
$ws.Cells.Item(22, 16).Value = 925

# Row 23
$ws.Cells.Item(23, 4).Value = 44455
$ws.Cells.Item(23, 10).Value = 1100
$ws.Cells.Item(23, 11).Value = 900
$ws.Cells.Item(23, 12).Value = 1000
$ws.Cells.Item(23, 13).Value = 950
$ws.Cells.Item(23, 16).Value = 950
